# Region XII_LMS.xlsx edit: add "Status as of July 4, 2025" column with a
# dropdown-validated list sourced from a new hidden "DropdownOptions" sheet,
# and tidy up a batch of previously-empty placeholder cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Clear out the stray empty placeholder cells that were left behind from a
#    prior export (they currently serialize as empty inline-string cells).
# ---------------------------------------------------------------------------
$emptyCells = @(
  "AU2","AV2","AW2","AX2","AY2","AZ2","BB2","BC2",
  "AV3","AW3","AX3","AY3","AZ3","BB3","BC3",
  "AU4","AW4","AX4","AY4","AZ4","BA4","BB4","BC4",
  "H5","AU5","AW5","AX5","AY5","AZ5","BA5","BB5","BC5",
  "AT6","AU6","AX6","AY6","AZ6","BA6","BB6","BC6",
  "AT7","AU7","AX7","AY7","AZ7","BA7","BB7","BC7",
  "O8","P8","S8","T8","U8","V8","W8","AT8","AU8","AV8","AW8","AX8","AY8","AZ8","BA8","BB8",
  "P9","AT9","AU9","AV9","AW9","AX9","AY9","AZ9","BA9","BB9",
  "P10","AT10","AU10","AV10","AW10","AX10","AY10","AZ10","BA10","BB10",
  "AT11","AU11","AV11","AW11","AX11","AY11","AZ11","BA11","BB11"
)
# NOTE: multi-area Range("A1,B2,...") unions only reliably clear the first
# area in this host, so clear each stray cell individually instead.
foreach ($addr in $emptyCells) {
  $ws.Range($addr).ClearContents()
}

# ---------------------------------------------------------------------------
# 2. Re-apply the date/time number format to the scheduling date columns so
#    they pick up a fresh, consolidated style entry.
# ---------------------------------------------------------------------------
$dateCells = @(
  "O9","S9","T9","U9","V9","W9",
  "O10","S10","T10","U10","V10","W10",
  "O11","S11","T11","U11","V11","W11"
)
foreach ($addr in $dateCells) {
  $c = $ws.Range($addr)
  $c.NumberFormat = "YYYY-MM-DD HH:MM:SS"
  $c.WrapText = $false
}

# ---------------------------------------------------------------------------
# 3. Add the hidden "DropdownOptions" sheet (construction-progress ranges)
#    right after Sheet1.
# ---------------------------------------------------------------------------
$dropdown = $wb.Worksheets.Add($null, $ws)
$dropdown.Name = "DropdownOptions"

$options = @(
  "0% - 10%: Foundation completed: Groundwork finished; no vertical structure yet.",
  "11% - 25%: Structure and rough-in started: Structural framing in progress; initial MEP rough-in.",
  "26% - 50%: Structure erected, partial roofing: Building shape defined; roof and systems advancing.",
  "51% - 75%: Exterior sealed, interior work underway: Enclosed structure; painting, flooring, and testing begin.",
  "76% - 90%: Final finishes and inspections: Systems tested; ",
  "91% - 99%: Final touches and punch list: Minor adjustments; final inspections and approvals.",
  "100% - Construction complete: Ready for handover and occupancy."
)
for ($i = 0; $i -lt $options.Length; $i++) {
  $dropdown.Cells.Item($i + 1, 1).Value = $options[$i]
}

$dropdown.Visible = $false

# ---------------------------------------------------------------------------
# 4. Add the new "Status as of July 4, 2025" header column on Sheet1 and wire
#    up the list-based data validation against the hidden sheet.
# ---------------------------------------------------------------------------
$ws.Range("BE1").Value = "Status as of July 4, 2025"

$statusRange = $ws.Range("BE2:BE11")
$statusRange.Validation.Add(3, 1, 1, "=DropdownOptions!`$A`$1:`$A`$7")
$statusRange.Validation.ShowInput = $false
$statusRange.Validation.ShowError = $false
$statusRange.Validation.InCellDropdown = $true

# Keep Sheet1 as the active tab (adding the new sheet shifts focus to it).
$ws.Activate()

Write-Output "edit complete"
